# devDocs.docx edit: merge split runs back into single runs (no visible
# text change) for a batch of paragraphs, and apply the one genuine text
# change - inserting "Google authentication" into the "Other
# authentication methods" bullet.

$d = $word.ActiveDocument
$dash = [char]0x2013

# Simple run-merge paragraphs: the replacement text is identical to the
# original text, which causes Word to coalesce the matched range back
# into a single run.
$d.Content.Find.Execute("Install node and all the dependencies including devDependencies.", $true, $false, $false, $false, $false, $true, 1, $false, "Install node and all the dependencies including devDependencies.", 2)

$d.Content.Find.Execute("Create an .env file inside the main project directory if not present already.(the content of this file will be emailed to professor seperately as it contains confidentail information).", $true, $false, $false, $false, $false, $true, 1, $false, "Create an .env file inside the main project directory if not present already.(the content of this file will be emailed to professor seperately as it contains confidentail information).", 2)

$d.Content.Find.Execute("Run command, once you are inside the project directory.", $true, $false, $false, $false, $false, $true, 1, $false, "Run command, once you are inside the project directory.", 2)

$d.Content.Find.Execute("Now you can check the project on the localhost:3000/", $true, $false, $false, $false, $false, $true, 1, $false, "Now you can check the project on the localhost:3000/", 2)

$d.Content.Find.Execute("Project includes various routes which perform typical functions :", $true, $false, $false, $false, $false, $true, 1, $false, "Project includes various routes which perform typical functions :", 2)

$d.Content.Find.Execute("Donate route", $true, $false, $false, $false, $false, $true, 1, $false, "Donate route", 2)

$d.Content.Find.Execute("This route check for the form validation and fetches data like country and states from another server. For validation - express-validator , creditcards, vanilla javascript is used. Sanatized data is inserted into database collection named - donate", $true, $false, $false, $false, $false, $true, 1, $false, "This route check for the form validation and fetches data like country and states from another server. For validation - express-validator , creditcards, vanilla javascript is used. Sanatized data is inserted into database collection named - donate", 2)

$d.Content.Find.Execute("This route displays the facebook profile of user.", $true, $false, $false, $false, $false, $true, 1, $false, "This route displays the facebook profile of user.", 2)

$d.Content.Find.Execute("Contact route", $true, $false, $false, $false, $false, $true, 1, $false, "Contact route", 2)

$d.Content.Find.Execute("This route sends feedback email after successfully filling out the contact form. It uses nodemailer for sending email and is associated with savewild email id for accomplishing aforementioned task. It stores contact queries inside the collection named " + $dash + " contact.", $true, $false, $false, $false, $false, $true, 1, $false, "This route sends feedback email after successfully filling out the contact form. It uses nodemailer for sending email and is associated with savewild email id for accomplishing aforementioned task. It stores contact queries inside the collection named " + $dash + " contact.", 2)

$d.Content.Find.Execute("Other routes", $true, $false, $false, $false, $false, $true, 1, $false, "Other routes", 2)

$d.Content.Find.Execute("What can be extended?", $true, $false, $false, $false, $false, $true, 1, $false, "What can be extended?", 2)

$d.Content.Find.Execute("The following upgradation can be made to the project:", $true, $false, $false, $false, $false, $true, 1, $false, "The following upgradation can be made to the project:", 2)

$d.Content.Find.Execute("Shopping cart with merchandise can be implemented.", $true, $false, $false, $false, $false, $true, 1, $false, "Shopping cart with merchandise can be implemented.", 2)

$d.Content.Find.Execute("Newsletter and reminder email can be send periodically to the users present inside the database", $true, $false, $false, $false, $false, $true, 1, $false, "Newsletter and reminder email can be send periodically to the users present inside the database", 2)

# Genuine content change (commit: "screencast and userDocs addedd") -
# the authentication bullet gains a mention of Google authentication.
$d.Content.Find.Execute("Other authentication methods can be used in the login route such as " + $dash + " twitter authentication , local authentication.", $true, $false, $false, $false, $false, $true, 1, $false, "Other authentication methods can be used in the login route such as " + $dash + " twitter authentication, Google authentication , local authentication.", 2)
